$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 (theta_se): B4 unchanged ("(0.63)" / nan placeholder), fill in the rest
$ws.Range("B4").Value = "(0.63)"
$ws.Range("C4").Value = "(0.52)"
$ws.Range("D4").Value = "(1.33)"
$ws.Range("E4").Value = "(0.09)"
$ws.Range("F4").Value = "(0.18)"
$ws.Range("G4").Value = "(0.46)"
$ws.Range("H4").Value = "(0.66)"
$ws.Range("I4").Value = "(0.86)"
$ws.Range("J4").Value = "(0.85)"

# Row 6 (lambda_se)
$ws.Range("B6").Value = "(0.6)"
$ws.Range("C6").Value = "(0.37)"
$ws.Range("D6").Value = "(1.22)"
$ws.Range("E6").Value = "(0.01)"
$ws.Range("F6").Value = "(0.01)"
$ws.Range("G6").Value = "(0.43)"
$ws.Range("H6").Value = "(0.19)"
$ws.Range("I6").Value = "(0.77)"
$ws.Range("J6").Value = "(0.0)"
